# Added titles for new Festival alts of Gold City and Symboli Rudolf
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("uma-title")

# Append the two new rows at the bottom of the uma-title table (A53:B54)
$ws.Range("A53").Value = "秋桜ダンツァトリーチェ"
$ws.Range("B53").Value = "Akizakura Danzatrice"
$ws.Range("A54").Value = "皓月の弓取り"
$ws.Range("B54").Value = "Archer of the White Moon"

# Grow the query table / ListObject so it covers the new rows
$lo = $ws.ListObjects("uma_title")
$lo.Resize($ws.Range("A1:B54"))

# Keep the hidden ExternalData_1 defined name in sync with the new extent
$wb.Names.Item("uma-title!ExternalData_1").RefersTo = "='uma-title'!`$A`$1:`$B`$54"

# Restore the selection used in the saved workbook
$ws.Range("D58").Select()
